$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Remove the two "CFM56-7B24" rows (rows 2 and 3); everything below shifts up by 2.
$ws.Rows("2:3").Delete()

# 2. Row 42 (P/N CFM56-7B26/3, S/N 802135) LLPLimit corrected to 33991,
#    now formatted like the EGTLimit column (#,##0.00 accounting style).
$c42 = $ws.Range("C42")
$c42.Value = 33991
$c42.NumberFormat = $ws.Range("D42").NumberFormat

# 3. Update the hidden _FilterDatabase defined name so its range reflects
#    the now-smaller table (39 data rows instead of 41).
foreach ($n in $wb.Names) {
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "=Sheet!`$A`$1:`$D`$39"
    }
}

# 4. Move/restore the active selection to B2 (matches author's last cursor spot).
[void]$ws.Range("B2").Select()
